# Generate Report for Handoff
#
# The handoff generation re-ran for the "35fcf230-f3f3-499f-8195-5edfd46dc5d4.md"
# file: its status flips from "In Translation" to "Ready for handoff" on every
# sheet, and the associated timestamps are bumped forward a few seconds.
# Because "Ready for handoff" is wider than "In Translation", Excel's autofit
# also grows the status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# Columns: A FileName, B PathAndName, C Extension, D PublishURL,
#          E zh-cn (status), F de-de (status), G Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 22:58:10"

# --- zh-cn sheet --------------------------------------------------------
# Column C is Status, column H is Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 22:58:05"

# --- de-de sheet --------------------------------------------------------
# Column C is Status, column H is Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 22:58:10"

# --- Re-fit the status columns that now hold the longer text -----------
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
